$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell-level updates describing the edit from the commit
# "graficos anuario y rraa DT": category-code relabeling (D->DE, H->HJ,
# L->LMN, R->RS, S->T) plus refreshed statistics for the affected rows
# across all ten year blocks (rows 5-161).
$updates = @(
    @{ Row = 5; Col = "A"; Value = "DE" },
    @{ Row = 8; Col = "A"; Value = "HJ" },
    @{ Row = 11; Col = "A"; Value = "LMN" },
    @{ Row = 15; Col = "A"; Value = "RS" },
    @{ Row = 15; Col = "C"; Value = 722868.3415958595 },
    @{ Row = 15; Col = "D"; Value = 71.247 },
    @{ Row = 15; Col = "E"; Value = 43.43334973716274 },
    @{ Row = 15; Col = "F"; Value = 32.635 },
    @{ Row = 15; Col = "G"; Value = 8.718999999999999 },
    @{ Row = 15; Col = "H"; Value = 57.611 },
    @{ Row = 15; Col = "I"; Value = 6.05 },
    @{ Row = 15; Col = "J"; Value = 5.623632311183862 },
    @{ Row = 15; Col = "K"; Value = 64.31 },
    @{ Row = 16; Col = "A"; Value = "T" },
    @{ Row = 16; Col = "C"; Value = 2574.58448395478 },
    @{ Row = 16; Col = "D"; Value = 44.375 },
    @{ Row = 16; Col = "E"; Value = 36.20396180932696 },
    @{ Row = 16; Col = "F"; Value = $null },
    @{ Row = 16; Col = "G"; Value = 68.04300000000001 },
    @{ Row = 16; Col = "H"; Value = 63.355 },
    @{ Row = 16; Col = "I"; Value = 54.543 },
    @{ Row = 16; Col = "J"; Value = 2.66911345069883 },
    @{ Row = 16; Col = "K"; Value = 71.92 },
    @{ Row = 17; Col = "E"; Value = 39.73855448822508 },
    @{ Row = 17; Col = "J"; Value = 5.690695256513068 },
    @{ Row = 21; Col = "A"; Value = "DE" },
    @{ Row = 24; Col = "A"; Value = "HJ" },
    @{ Row = 27; Col = "A"; Value = "LMN" },
    @{ Row = 31; Col = "A"; Value = "RS" },
    @{ Row = 31; Col = "C"; Value = 767966.6091438441 },
    @{ Row = 31; Col = "D"; Value = 72.84999999999999 },
    @{ Row = 31; Col = "E"; Value = 44.58761596695263 },
    @{ Row = 31; Col = "F"; Value = 29.472 },
    @{ Row = 31; Col = "G"; Value = 8.922000000000001 },
    @{ Row = 31; Col = "H"; Value = 60.645 },
    @{ Row = 31; Col = "I"; Value = 6.148 },
    @{ Row = 31; Col = "J"; Value = 5.61223235570472 },
    @{ Row = 31; Col = "K"; Value = 67.53100000000001 },
    @{ Row = 32; Col = "A"; Value = "T" },
    @{ Row = 32; Col = "C"; Value = 2121.714470077975 },
    @{ Row = 32; Col = "D"; Value = 63.785 },
    @{ Row = 32; Col = "E"; Value = 46.99070136263129 },
    @{ Row = 32; Col = "F"; Value = $null },
    @{ Row = 32; Col = "G"; Value = 60.212 },
    @{ Row = 32; Col = "H"; Value = 97.706 },
    @{ Row = 32; Col = "I"; Value = 59.462 },
    @{ Row = 32; Col = "J"; Value = 7.368114738954834 },
    @{ Row = 32; Col = "K"; Value = 100 },
    @{ Row = 33; Col = "E"; Value = 40.62391110297622 },
    @{ Row = 33; Col = "J"; Value = 5.836243566245659 },
    @{ Row = 37; Col = "A"; Value = "DE" },
    @{ Row = 40; Col = "A"; Value = "HJ" },
    @{ Row = 43; Col = "A"; Value = "LMN" },
    @{ Row = 47; Col = "A"; Value = "RS" },
    @{ Row = 47; Col = "C"; Value = 729429.3514179224 },
    @{ Row = 47; Col = "D"; Value = 75.68300000000001 },
    @{ Row = 47; Col = "E"; Value = 44.80379179463692 },
    @{ Row = 47; Col = "F"; Value = 30.958 },
    @{ Row = 47; Col = "G"; Value = 7.645 },
    @{ Row = 47; Col = "H"; Value = 59.602 },
    @{ Row = 47; Col = "I"; Value = 6.627 },
    @{ Row = 47; Col = "J"; Value = 5.54562362824567 },
    @{ Row = 47; Col = "K"; Value = 64.553 },
    @{ Row = 48; Col = "A"; Value = "T" },
    @{ Row = 48; Col = "C"; Value = 1916.540108464908 },
    @{ Row = 48; Col = "D"; Value = 65.26600000000001 },
    @{ Row = 48; Col = "E"; Value = 41.92468395004806 },
    @{ Row = 48; Col = "F"; Value = $null },
    @{ Row = 48; Col = "G"; Value = 46.469 },
    @{ Row = 48; Col = "H"; Value = 88.354 },
    @{ Row = 48; Col = "I"; Value = 62.55 },
    @{ Row = 48; Col = "J"; Value = 8.077464662531003 },
    @{ Row = 48; Col = "K"; Value = 81.44 },
    @{ Row = 49; Col = "E"; Value = 40.38993185012123 },
    @{ Row = 49; Col = "J"; Value = 5.891564841054532 },
    @{ Row = 53; Col = "A"; Value = "DE" },
    @{ Row = 56; Col = "A"; Value = "HJ" },
    @{ Row = 59; Col = "A"; Value = "LMN" },
    @{ Row = 63; Col = "A"; Value = "RS" },
    @{ Row = 63; Col = "C"; Value = 740682.9000341403 },
    @{ Row = 63; Col = "D"; Value = 72.617 },
    @{ Row = 63; Col = "E"; Value = 45.11333319933395 },
    @{ Row = 63; Col = "F"; Value = 30.539 },
    @{ Row = 63; Col = "G"; Value = 11.146 },
    @{ Row = 63; Col = "H"; Value = 58.919 },
    @{ Row = 63; Col = "I"; Value = 6.187 },
    @{ Row = 63; Col = "J"; Value = 5.944428224247003 },
    @{ Row = 63; Col = "K"; Value = 66.864 },
    @{ Row = 64; Col = "A"; Value = "T" },
    @{ Row = 64; Col = "C"; Value = 3276.909896187492 },
    @{ Row = 64; Col = "D"; Value = 57.478 },
    @{ Row = 64; Col = "E"; Value = 40.34438217287427 },
    @{ Row = 64; Col = "F"; Value = $null },
    @{ Row = 64; Col = "G"; Value = 70.473 },
    @{ Row = 64; Col = "H"; Value = 56.125 },
    @{ Row = 64; Col = "I"; Value = 56.541 },
    @{ Row = 64; Col = "J"; Value = 6.385025095665147 },
    @{ Row = 64; Col = "K"; Value = 87.27200000000001 },
    @{ Row = 65; Col = "E"; Value = 40.90595833638884 },
    @{ Row = 65; Col = "J"; Value = 6.103717189424323 },
    @{ Row = 69; Col = "A"; Value = "DE" },
    @{ Row = 72; Col = "A"; Value = "HJ" },
    @{ Row = 75; Col = "A"; Value = "LMN" },
    @{ Row = 79; Col = "A"; Value = "RS" },
    @{ Row = 79; Col = "C"; Value = 731477.9944147074 },
    @{ Row = 79; Col = "D"; Value = 73.06 },
    @{ Row = 79; Col = "E"; Value = 45.57726092856868 },
    @{ Row = 79; Col = "F"; Value = 29.165 },
    @{ Row = 79; Col = "G"; Value = 9.872 },
    @{ Row = 79; Col = "H"; Value = 60.823 },
    @{ Row = 79; Col = "I"; Value = 7.911 },
    @{ Row = 79; Col = "J"; Value = 5.94333646036294 },
    @{ Row = 79; Col = "K"; Value = 65.914 },
    @{ Row = 80; Col = "A"; Value = "T" },
    @{ Row = 80; Col = "C"; Value = 1034.092456731289 },
    @{ Row = 80; Col = "D"; Value = 84.908 },
    @{ Row = 80; Col = "E"; Value = 47.50875731870244 },
    @{ Row = 80; Col = "F"; Value = $null },
    @{ Row = 80; Col = "G"; Value = 29.441 },
    @{ Row = 80; Col = "H"; Value = 84.854 },
    @{ Row = 80; Col = "I"; Value = 100 },
    @{ Row = 80; Col = "J"; Value = 8.366572959205735 },
    @{ Row = 80; Col = "K"; Value = 100 },
    @{ Row = 81; Col = "E"; Value = 41.3077207655281 },
    @{ Row = 81; Col = "J"; Value = 6.082950159367544 },
    @{ Row = 85; Col = "A"; Value = "DE" },
    @{ Row = 88; Col = "A"; Value = "HJ" },
    @{ Row = 91; Col = "A"; Value = "LMN" },
    @{ Row = 95; Col = "A"; Value = "RS" },
    @{ Row = 95; Col = "C"; Value = 668564.3980730547 },
    @{ Row = 95; Col = "D"; Value = 74.07899999999999 },
    @{ Row = 95; Col = "E"; Value = 46.7295888940445 },
    @{ Row = 95; Col = "F"; Value = 25.858 },
    @{ Row = 95; Col = "G"; Value = 9.455 },
    @{ Row = 95; Col = "H"; Value = 64.476 },
    @{ Row = 95; Col = "I"; Value = 7.126 },
    @{ Row = 95; Col = "J"; Value = 5.808815176889641 },
    @{ Row = 95; Col = "K"; Value = 70.383 },
    @{ Row = 96; Col = "A"; Value = "T" },
    @{ Row = 96; Col = "C"; Value = 2422.251365622006 },
    @{ Row = 96; Col = "D"; Value = 42.742 },
    @{ Row = 96; Col = "E"; Value = 45.69400371699194 },
    @{ Row = 96; Col = "F"; Value = $null },
    @{ Row = 96; Col = "G"; Value = 55.011 },
    @{ Row = 96; Col = "H"; Value = 93.851 },
    @{ Row = 96; Col = "I"; Value = 25.892 },
    @{ Row = 96; Col = "J"; Value = 6.838084022341369 },
    @{ Row = 96; Col = "K"; Value = 93.851 },
    @{ Row = 97; Col = "E"; Value = 41.55580109936026 },
    @{ Row = 97; Col = "J"; Value = 5.804348038843218 },
    @{ Row = 101; Col = "A"; Value = "DE" },
    @{ Row = 104; Col = "A"; Value = "HJ" },
    @{ Row = 107; Col = "A"; Value = "LMN" },
    @{ Row = 111; Col = "A"; Value = "RS" },
    @{ Row = 111; Col = "C"; Value = 690995.569793692 },
    @{ Row = 111; Col = "D"; Value = 71.496 },
    @{ Row = 111; Col = "E"; Value = 45.1172180464615 },
    @{ Row = 111; Col = "F"; Value = 21.855 },
    @{ Row = 111; Col = "G"; Value = 10.519 },
    @{ Row = 111; Col = "H"; Value = 67.79900000000001 },
    @{ Row = 111; Col = "I"; Value = 6.951 },
    @{ Row = 111; Col = "J"; Value = 6.414328450327223 },
    @{ Row = 111; Col = "K"; Value = 69.625 },
    @{ Row = 112; Col = "A"; Value = "T" },
    @{ Row = 112; Col = "C"; Value = 1083.007952034158 },
    @{ Row = 112; Col = "D"; Value = 98.705 },
    @{ Row = 112; Col = "E"; Value = 40.44142440450224 },
    @{ Row = 112; Col = "F"; Value = $null },
    @{ Row = 112; Col = "G"; Value = 77.248 },
    @{ Row = 112; Col = "H"; Value = 43.559 },
    @{ Row = 112; Col = "I"; Value = 41.026 },
    @{ Row = 112; Col = "J"; Value = 8.639761769387789 },
    @{ Row = 112; Col = "K"; Value = 100 },
    @{ Row = 113; Col = "E"; Value = 41.49119268525506 },
    @{ Row = 113; Col = "J"; Value = 6.315529348701673 },
    @{ Row = 117; Col = "A"; Value = "DE" },
    @{ Row = 120; Col = "A"; Value = "HJ" },
    @{ Row = 123; Col = "A"; Value = "LMN" },
    @{ Row = 127; Col = "A"; Value = "RS" },
    @{ Row = 127; Col = "C"; Value = 730217.0869155917 },
    @{ Row = 127; Col = "D"; Value = 71.508 },
    @{ Row = 127; Col = "E"; Value = 44.99873051821902 },
    @{ Row = 127; Col = "F"; Value = 21.459 },
    @{ Row = 127; Col = "G"; Value = 10.969 },
    @{ Row = 127; Col = "H"; Value = 66.64 },
    @{ Row = 127; Col = "I"; Value = 9.257 },
    @{ Row = 127; Col = "J"; Value = 5.946262875100936 },
    @{ Row = 127; Col = "K"; Value = 68.917 },
    @{ Row = 128; Col = "A"; Value = "T" },
    @{ Row = 128; Col = "C"; Value = 4369.181973960196 },
    @{ Row = 128; Col = "D"; Value = 58.422 },
    @{ Row = 128; Col = "E"; Value = 43.08693594143686 },
    @{ Row = 128; Col = "F"; Value = $null },
    @{ Row = 128; Col = "G"; Value = 41.334 },
    @{ Row = 128; Col = "H"; Value = 90.18600000000001 },
    @{ Row = 128; Col = "I"; Value = 50.694 },
    @{ Row = 128; Col = "J"; Value = 3.175709624834065 },
    @{ Row = 128; Col = "K"; Value = 77.166 },
    @{ Row = 129; Col = "E"; Value = 41.94639353121017 },
    @{ Row = 129; Col = "J"; Value = 5.937164596229016 },
    @{ Row = 133; Col = "A"; Value = "DE" },
    @{ Row = 136; Col = "A"; Value = "HJ" },
    @{ Row = 139; Col = "A"; Value = "LMN" },
    @{ Row = 143; Col = "A"; Value = "RS" },
    @{ Row = 143; Col = "C"; Value = 706194.5280358862 },
    @{ Row = 143; Col = "D"; Value = 69.75700000000001 },
    @{ Row = 143; Col = "E"; Value = 45.94104659275887 },
    @{ Row = 143; Col = "F"; Value = 23.649 },
    @{ Row = 143; Col = "G"; Value = 13.212 },
    @{ Row = 143; Col = "H"; Value = 65.771 },
    @{ Row = 143; Col = "I"; Value = 8.315 },
    @{ Row = 143; Col = "J"; Value = 6.084848949609335 },
    @{ Row = 143; Col = "K"; Value = 67.327 },
    @{ Row = 144; Col = "A"; Value = "T" },
    @{ Row = 144; Col = "C"; Value = 1110.94803144716 },
    @{ Row = 144; Col = "D"; Value = 76.051 },
    @{ Row = 144; Col = "E"; Value = 45.37439076018008 },
    @{ Row = 144; Col = "F"; Value = $null },
    @{ Row = 144; Col = "G"; Value = 69.133 },
    @{ Row = 144; Col = "H"; Value = 100 },
    @{ Row = 144; Col = "I"; Value = 60.176 },
    @{ Row = 144; Col = "J"; Value = 7.982231415218051 },
    @{ Row = 144; Col = "K"; Value = 100 },
    @{ Row = 145; Col = "E"; Value = 42.31281225420921 },
    @{ Row = 145; Col = "J"; Value = 6.122303416440997 },
    @{ Row = 149; Col = "A"; Value = "DE" },
    @{ Row = 152; Col = "A"; Value = "HJ" },
    @{ Row = 155; Col = "A"; Value = "LMN" },
    @{ Row = 159; Col = "A"; Value = "RS" },
    @{ Row = 159; Col = "C"; Value = 705715.5041364761 },
    @{ Row = 159; Col = "D"; Value = 69.54000000000001 },
    @{ Row = 159; Col = "E"; Value = 46.32042738309166 },
    @{ Row = 159; Col = "F"; Value = 25.375 },
    @{ Row = 159; Col = "G"; Value = 12.13 },
    @{ Row = 159; Col = "H"; Value = 67.122 },
    @{ Row = 159; Col = "I"; Value = 7.78 },
    @{ Row = 159; Col = "J"; Value = 5.788153731572223 },
    @{ Row = 159; Col = "K"; Value = 62.441 },
    @{ Row = 160; Col = "A"; Value = "T" },
    @{ Row = 160; Col = "C"; Value = 1002.225371042404 },
    @{ Row = 160; Col = "D"; Value = 73.11499999999999 },
    @{ Row = 160; Col = "E"; Value = 48.94318887129524 },
    @{ Row = 160; Col = "F"; Value = $null },
    @{ Row = 160; Col = "G"; Value = 35.821 },
    @{ Row = 160; Col = "H"; Value = 8.936 },
    @{ Row = 160; Col = "I"; Value = 35.821 },
    @{ Row = 160; Col = "J"; Value = 0.8065636888976946 },
    @{ Row = 160; Col = "K"; Value = 100 },
    @{ Row = 161; Col = "E"; Value = 42.71633487421943 },
    @{ Row = 161; Col = "J"; Value = 5.454940117397035 }
)

foreach ($item in $updates) {
    $addr = "$($item.Col)$($item.Row)"
    if ($item.Value -eq $null) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $item.Value
    }
}
